# Creation allocation matrix for dMRWIO + Results_region
#
# For every year sheet (2000..2100), the 3-row "Sector" block (rows 5-7,
# columns C..G) lists Onshore wind plants / Offshore wind plants /
# Photovoltaic plants with their corresponding value in column E.
# The edit re-orders that block so "Photovoltaic plants" comes first,
# followed by "Onshore wind plants" then "Offshore wind plants" -
# i.e. row 7 is moved above row 5, pushing the other two rows down by one.
#
# Concretely: new row5 = old row7, new row6 = old row5, new row7 = old row6.
# Columns D, F, G are always 0 and stay untouched; column A/B are merged
# across rows 5:7 and also stay untouched.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $c5 = $ws.Range("C5").Value2
    $c6 = $ws.Range("C6").Value2
    $c7 = $ws.Range("C7").Value2

    $e5 = $ws.Range("E5").Value2
    $e6 = $ws.Range("E6").Value2
    $e7 = $ws.Range("E7").Value2

    # Rotate the 3-row block down by one (old row7 -> new row5, etc.)
    $ws.Range("C5").Value2 = $c7
    $ws.Range("C6").Value2 = $c5
    $ws.Range("C7").Value2 = $c6

    $ws.Range("E5").Value2 = $e7
    $ws.Range("E6").Value2 = $e5
    $ws.Range("E7").Value2 = $e6
}
